# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# for the rows whose market data changed. A leading apostrophe forces the
# value to be stored as text (matching the original inline-string cells,
# since some prices like "213.21" would otherwise be auto-parsed as a
# number by Excel); ClearFormats() afterwards drops the quote-prefix style
# flag so the cell keeps using the sheet default style, unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.659.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +1.36%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.634.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.92%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D5").Value = "'213.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.57%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.501"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +3.61%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.253"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.26%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +1.37%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'19.21"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +2.15%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +3.52%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.862.66"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'1.620.35"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.11%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'  +2.24%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  +1.54%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'26.651.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.33%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'63.39"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.87%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'  +2.24%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'219.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +8.94%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  +0.03%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  +0.37%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'9.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +1.42%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'6.21"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +2.66%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  +2.20%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'148.73"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +2.98%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +1.41%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +5.68%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'15.53"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +2.49%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -2.11%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  +0.06%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'3.30"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +3.82%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  +2.65%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +0.98%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  -0.05%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'1.213.82"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +2.92%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +5.52%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  +0.84%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  +0.08%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.505"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +1.87%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  -1.17%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'5.42"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +1.48%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.794"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.72%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'1.770.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.85%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'93.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +0.60%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'  +1.14%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'54.76"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +1.95%  "
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'  +0.94%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'7.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +5.36%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  +0.42%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  +0.20%  "
$ws.Range("E51").ClearFormats()
